$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '21.648.93'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.90%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.533.81'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.55%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '288.31'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.28%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3936'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.26%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3158'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.73%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.43'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07160'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.73%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -7.26%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.10%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.652'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.54'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -4.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.584'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.539.73'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001089'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.55%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06605'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '83.51'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.00%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.110'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -4.73%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.42'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -3.62%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.78'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -6.22%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.348'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.64%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '21.651.44'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.90%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.350'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -7.73%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '148.21'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.53%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.31'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.98%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.830'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.59%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.710.46'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.22%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '117.06'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.936'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.81%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9455'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -14.99%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08144'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.43%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '8.482'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -8.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.128'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.70%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06000'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.57%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02203'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -4.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.443'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -14.37%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.57%  '
$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '10.99'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.78%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.172'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -4.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9997'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5747'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -3.68%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.05'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.705'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.35%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5477'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -4.87%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.163'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.26%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.869'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.43%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '116.13'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.66%  '
